$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("ref OEM" shifts from B to C,
# "Nouveau prix d'achat" shifts from C to D), then set the new header.
$ws.Range("B1").EntireColumn.Insert()
$ws.Range("B1").Value = "libellé"

# Reset the selection to the default (A1) so the sheet is saved without an
# explicit <selection> override.
$ws.Range("A1").Select()
